# Add a new drug record (grazoprevir / GZR / MK-5172) as row 10 of Sheet1,
# following the same pattern as the existing rows (id, abbreviation,
# category, producer, researchCode).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: the existing "NS3/4A protease inhibitors" category string in the
# workbook's shared-string table uses a non-breaking space (U+00A0) between
# "NS3/4A" and "protease" -- reproduce that exactly so Excel reuses the
# existing shared string instead of creating a near-duplicate one.
$nbspCategory = "NS3/4A" + [char]0x00A0 + "protease inhibitors"

$ws.Range("A10").Value = "grazoprevir"
$ws.Range("B10").Value = "GZR"
$ws.Range("C10").Value = $nbspCategory
$ws.Range("D10").Value = "Merck"
$ws.Range("E10").Value = "MK-5172"

# Match the formatting already used in column B for the preceding two rows
# (B8, B9) so the new cell picks up the same style.
$ws.Range("B10").Style = $ws.Range("B9").Style
